$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Force text (NumberFormat "@") before assigning values that could
# otherwise be auto-coerced to numbers by Excel (e.g. "1.005", "314.45"),
# then reset the style back to Normal so no stray formatting is left
# behind (matches the original unstyled inline-string cells).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '27.248.48'
Set-TextValue 'E2' '  -1.47%  '
Set-TextValue 'D3' '1.821.71'
Set-TextValue 'E3' '  -1.93%  '
Set-TextValue 'E4' '  -1.44%  '
Set-TextValue 'D5' '314.45'
Set-TextValue 'E5' '  -1.83%  '
Set-TextValue 'D6' '1.005'
Set-TextValue 'E6' '  -1.24%  '
Set-TextValue 'D7' '0.4266'
Set-TextValue 'E7' '  -2.25%  '
Set-TextValue 'D8' '0.3669'
Set-TextValue 'E8' '  -3.17%  '
Set-TextValue 'D9' '45.99'
Set-TextValue 'E9' '  -1.59%  '
Set-TextValue 'D10' '0.07215'
Set-TextValue 'E10' '  -2.77%  '
Set-TextValue 'D11' '0.8603'
Set-TextValue 'E11' '  -2.62%  '
Set-TextValue 'D12' '20.97'
Set-TextValue 'E12' '  -3.02%  '
Set-TextValue 'D13' '1.842.50'
Set-TextValue 'E13' '  -1.00%  '
Set-TextValue 'D14' '6.650'
Set-TextValue 'E14' '  -1.25%  '
Set-TextValue 'D15' '0.07102'
Set-TextValue 'E15' '  +0.10%  '
Set-TextValue 'D16' '5.301'
Set-TextValue 'E16' '  -3.44%  '
Set-TextValue 'D17' '87.90'
Set-TextValue 'E17' '  +1.35%  '
Set-TextValue 'D18' '1.007'
Set-TextValue 'E18' '  -1.51%  '
Set-TextValue 'D19' '0.000008857'
Set-TextValue 'E19' '  -2.29%  '
Set-TextValue 'E20' '  -1.25%  '
Set-TextValue 'D21' '15.03'
Set-TextValue 'E21' '  -2.73%  '
Set-TextValue 'D22' '27.264.77'
Set-TextValue 'E22' '  -1.43%  '
Set-TextValue 'D23' '5.126'
Set-TextValue 'E23' '  -2.97%  '
Set-TextValue 'D24' '10.88'
Set-TextValue 'E24' '  -2.43%  '
Set-TextValue 'D25' '2.053.89'
Set-TextValue 'E25' '  -2.31%  '
Set-TextValue 'E26' '  -1.51%  '
Set-TextValue 'D27' '153.04'
Set-TextValue 'E27' '  -2.56%  '
Set-TextValue 'D28' '18.33'
Set-TextValue 'E28' '  -1.96%  '
Set-TextValue 'D29' '2.102'
Set-TextValue 'E29' '  +5.23%  '
Set-TextValue 'D30' '5.222'
Set-TextValue 'E30' '  -2.57%  '
Set-TextValue 'D31' '116.19'
Set-TextValue 'E31' '  -3.51%  '
Set-TextValue 'D32' '0.08860'
Set-TextValue 'E32' '  -2.13%  '
Set-TextValue 'B33' 'ImmutableX'
Set-TextValue 'C33' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D33' '0.7611'
Set-TextValue 'E33' '  -1.04%  '
Set-TextValue 'B34' 'ARBITRUM'
Set-TextValue 'C34' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D34' '1.192'
Set-TextValue 'E34' '  -2.12%  '
Set-TextValue 'D35' '4.449'
Set-TextValue 'E35' '  -2.32%  '
Set-TextValue 'D36' '2.795'
Set-TextValue 'E36' '  -7.68%  '
Set-TextValue 'D37' '1.005'
Set-TextValue 'E37' '  -1.43%  '
Set-TextValue 'E38' '  -2.32%  '
Set-TextValue 'D39' '0.01955'
Set-TextValue 'E39' '  -1.04%  '
Set-TextValue 'D40' '0.05251'
Set-TextValue 'E40' '  -0.85%  '
Set-TextValue 'D41' '2.900'
Set-TextValue 'E41' '  +1.22%  '
Set-TextValue 'D42' '7.046'
Set-TextValue 'E42' '  +1.34%  '
Set-TextValue 'D43' '0.1676'
Set-TextValue 'E43' '  -0.20%  '
Set-TextValue 'D44' '0.5013'
Set-TextValue 'E44' '  -3.45%  '
Set-TextValue 'D45' '8.614'
Set-TextValue 'E45' '  -0.81%  '
Set-TextValue 'D46' '10.53'
Set-TextValue 'E46' '  -2.36%  '
Set-TextValue 'D47' '106.45'
Set-TextValue 'E47' '  -3.23%  '
Set-TextValue 'D48' '0.4685'
Set-TextValue 'E48' '  -0.82%  '
Set-TextValue 'D49' '1.005'
Set-TextValue 'E49' '  -1.34%  '
Set-TextValue 'E50' '  -1.63%  '
Set-TextValue 'D51' '1.657'
Set-TextValue 'E51' '  -3.18%  '
